$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna5"
$ws.Range("C2").Value = "Epha5"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05800433333333333
$ws.Range("H2").Value = 0.174013
$ws.Range("I2").Value = 0.02087975181349295
$ws.Range("J2").Value = 0.02087975181349295
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.008504666666666666
$ws.Range("N2").Value = 0.025514
$ws.Range("O2").Value = 0.3779124020559003
$ws.Range("P2").Value = 0.3779124020559003
$ws.Range("Q2").Value = 0.0004933075202222222
$ws.Range("R2").Value = 0.004439767682
$ws.Range("S2").Value = 0.007890717162168161
$ws.Range("T2").Value = 0.007890717162168163

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efna5"
$ws.Range("C3").Value = "Epha5"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.05800433333333333
$ws.Range("H3").Value = 0.174013
$ws.Range("I3").Value = 0.02087975181349295
$ws.Range("J3").Value = 0.02087975181349295
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01399966666666667
$ws.Range("N3").Value = 0.041999
$ws.Range("O3").Value = 0.6220875979440996
$ws.Range("P3").Value = 0.6220875979440996
$ws.Range("Q3").Value = 0.0008120413318888889
$ws.Range("R3").Value = 0.007308371987
$ws.Range("S3").Value = 0.01298903465132479
$ws.Range("T3").Value = 0.01298903465132479

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Efna5"
$ws.Range("C4").Value = "Epha5"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.666083666666667
$ws.Range("H4").Value = 4.998251
$ws.Range("I4").Value = 0.5997381826733804
$ws.Range("J4").Value = 0.5997381826733805
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.008504666666666666
$ws.Range("N4").Value = 0.025514
$ws.Range("O4").Value = 0.3779124020559003
$ws.Range("P4").Value = 0.3779124020559003
$ws.Range("Q4").Value = 0.01416948622377777
$ws.Range("R4").Value = 0.127525376014
$ws.Range("S4").Value = 0.2266484972187375
$ws.Range("T4").Value = 0.2266484972187376

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efna5"
$ws.Range("C5").Value = "Epha5"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.666083666666667
$ws.Range("H5").Value = 4.998251
$ws.Range("I5").Value = 0.5997381826733804
$ws.Range("J5").Value = 0.5997381826733805
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01399966666666667
$ws.Range("N5").Value = 0.041999
$ws.Range("O5").Value = 0.6220875979440996
$ws.Range("P5").Value = 0.6220875979440996
$ws.Range("Q5").Value = 0.02332461597211111
$ws.Range("R5").Value = 0.209921543749
$ws.Range("S5").Value = 0.3730896854546428
$ws.Range("T5").Value = 0.3730896854546429

$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Efna5"
$ws.Range("C6").Value = "Epha5"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.053930333333333
$ws.Range("H6").Value = 3.161791
$ws.Range("I6").Value = 0.3793820655131266
$ws.Range("J6").Value = 0.3793820655131266
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.008504666666666666
$ws.Range("N6").Value = 0.025514
$ws.Range("O6").Value = 0.3779124020559003
$ws.Range("P6").Value = 0.3779124020559003
$ws.Range("Q6").Value = 0.008963326174888889
$ws.Range("R6").Value = 0.080669935574
$ws.Range("S6").Value = 0.1433731876749946
$ws.Range("T6").Value = 0.1433731876749946

$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Efna5"
$ws.Range("C7").Value = "Epha5"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.053930333333333
$ws.Range("H7").Value = 3.161791
$ws.Range("I7").Value = 0.3793820655131266
$ws.Range("J7").Value = 0.3793820655131266
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.01399966666666667
$ws.Range("N7").Value = 0.041999
$ws.Range("O7").Value = 0.6220875979440996
$ws.Range("P7").Value = 0.6220875979440996
$ws.Range("Q7").Value = 0.01475467335655556
$ws.Range("R7").Value = 0.132792060209
$ws.Range("S7").Value = 0.236008877838132
$ws.Range("T7").Value = 0.236008877838132
